$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for each record.
# Update the date value from 45179 (2023-09-10) to 45180 (2023-09-11)
# for every data row (rows 2 through 261).
$ws.Range("C2:C261").Value = 45180
